# feat: make first comment..
#
# Adds the first "Prologue/firstComment" localization row (row 8) to the
# Languages sheet: ID/Name:key, EN and RU text, widens/re-splits column C
# so the new (longer) EN text is readable, grows row 8's height to fit the
# wrapped text, and moves the active selection to C9 (just below the new
# row), matching the authored change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New localization entry in row 8 -------------------------------------
# A8/B8 share the same key string, like the existing MainMenu/* rows above.
$ws.Range("A8").Value = "Prologue/firstComment"
$ws.Range("B8").Value = "Prologue/firstComment"
$ws.Range("C8").Value = "If you dream of this ..."
$ws.Range("D8").Value = "Приснится же такое…"

# Match the formatting used by the other data rows (5-7): centered Arial.
$ws.Range("B8").HorizontalAlignment = -4108
$ws.Range("B8").VerticalAlignment = -4108
$ws.Range("B8").Font.Name = "Arial"
$ws.Range("B8").Font.Size = 11

# Row grows taller to fit the wrapped English comment text.
$ws.Rows.Item(8).RowHeight = 26.85

# --- Column layout ---------------------------------------------------------
# Column C (EN) gets its own, wider column so the new comment text fits;
# columns D+ keep the previous default width.
$ws.Columns.Item(3).ColumnWidth = 20.17

# --- Selection ---------------------------------------------------------
# Active cell moves to C9, right under the freshly-added row.
$ws.Range("C9").Select() | Out-Null
